# Cryptos list update - Sun Mar 17 09:16:32 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.576.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.558.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.27%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.553.82"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.668"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.147"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.67%  "
$ws.Range("E12").Value = "  -8.32%  "
$ws.Range("E13").Value = "  -11.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.124.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.560.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.31%  "
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "66.491.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.02%  "
$ws.Range("E21").Value = "  -8.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.69%  "
$ws.Range("E23").Value = "  -8.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.38%  "
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.35%  "
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.57%  "
$ws.Range("E30").Value = "  -9.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "625.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("E35").Value = "  -7.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.402"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0744"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.97%  "
$ws.Range("E41").Value = "  -5.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.079.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("E46").Value = "  -9.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.06%  "

# Row 47 <-> Row 48 content swap (Stellar / ApeXProtocol)
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.129"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.36%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.14%  "

# Row 50 <-> Row 51 content swap (Stacks / THORChain)
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.74%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.58%  "
